# Apply the latest cryptocurrency price/volume snapshot to the "cryptos" worksheet.
# Generated to match the upstream GitHub Actions data refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.282.89"
$ws.Range("E2").Value = "  +5.34%  "
$ws.Range("D3").Value = "2.032.52"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.29"
$ws.Range("E5").Value = "  +5.14%  "
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.24"
$ws.Range("E7").Value = "  +15.30%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("E9").Value = "  +6.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.82"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.904"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.91"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("D15").Value = "2.337.73"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +7.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.40"
$ws.Range("E17").Value = "  +21.56%  "
$ws.Range("D18").Value = "1.976.88"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "37.265.43"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.12"
$ws.Range("E20").Value = "  +5.32%  "
$ws.Range("D21").Value = "0.0₃0869"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.30"
$ws.Range("E22").Value = "  +8.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.74"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.75"
$ws.Range("E24").Value = "  +23.94%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.47"
$ws.Range("E27").Value = "  +5.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.29"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.75"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.16"
$ws.Range("E31").Value = "  +10.24%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.113"
$ws.Range("E32").Value = "  +27.28%  "
$ws.Range("E33").Value = "  +7.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  +12.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0612"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +13.46%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.95"
$ws.Range("E39").Value = "  +24.72%  "
$ws.Range("E40").Value = "  +20.71%  "
$ws.Range("E41").Value = "  +6.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0218"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  +6.74%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").Value = "  +22.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.02"
$ws.Range("E46").Value = "  +10.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.75"
$ws.Range("E47").Value = "  +11.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.68"
$ws.Range("E48").Value = "  +7.30%  "
$ws.Range("D49").Value = "1.415.46"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.16"
$ws.Range("E51").Value = "  +6.79%  "

Write-Output "cryptos worksheet updated"
